# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-19, replacing the previous
# Strike# derived values with the regenerated K values.
$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
